$d = $word.ActiveDocument

function Set-ParagraphRuns {
    param(
        $Paragraph,
        [System.Collections.IEnumerable]$Runs
    )
    $r = $Paragraph.Range.Duplicate
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = ""
    $r.Collapse(1) | Out-Null

    foreach ($run in $Runs) {
        $kind = $run[0]
        if ($kind -eq "text") {
            $text = $run[1]
            $style = $run[2]
            $r.InsertAfter($text)
            if ($style) {
                $r.Style = $style
            }
            $r.Collapse(0) | Out-Null
        } elseif ($kind -eq "br") {
            $r.InsertAfter([char]11)
            $r.Collapse(0) | Out-Null
        }
    }
}

# ------------------------------------------------------------------
# 1. Title paragraph: "Welch's ANOVA and Kruskal-Wallis Test for
#    GSIand PCI" -> "ANOVA and Kruskal-Wallis Test for GSI and PCI"
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
Set-ParagraphRuns $titlePara @(
    ,("text", "ANOVA", $null)
    ,("text", " ", $null)
    ,("text", "and", $null)
    ,("text", " ", $null)
    ,("text", "Kruskal-Wallis", $null)
    ,("text", " ", $null)
    ,("text", "Test", $null)
    ,("text", " ", $null)
    ,("text", "for", $null)
    ,("text", " ", $null)
    ,("text", "GSI", $null)
    ,("text", " ", $null)
    ,("text", "and", $null)
    ,("text", " ", $null)
    ,("text", "PCI", $null)
)

# ------------------------------------------------------------------
# 2. Heading2: "Welch's ANOVA (Handles Unequal Variances)" ->
#    "Standard One-Way ANOVA"
# ------------------------------------------------------------------
$headingPara = $d.Paragraphs.Item(17)
Set-ParagraphRuns $headingPara @(
    ,("text", "Standard One-Way ANOVA", $null)
)

# ------------------------------------------------------------------
# 3. GSI aov() call paragraph
# ------------------------------------------------------------------
$gsiCallPara = $d.Paragraphs.Item(18)
Set-ParagraphRuns $gsiCallPara @(
    ,("text", "gsi_anova ", "NormalTok")
    ,("text", "<-", "OtherTok")
    ,("text", " ", "NormalTok")
    ,("text", "aov", "FunctionTok")
    ,("text", "(GSI ", "NormalTok")
    ,("text", "~", "SpecialCharTok")
    ,("text", " Size_Class, ", "NormalTok")
    ,("text", "data =", "AttributeTok")
    ,("text", " cots_data)", "NormalTok")
    ,("br", $null, $null)
    ,("text", "summary", "FunctionTok")
    ,("text", "(gsi_anova)", "NormalTok")
)

# ------------------------------------------------------------------
# 4. GSI summary() output paragraph
# ------------------------------------------------------------------
$gsiOutPara = $d.Paragraphs.Item(19)
Set-ParagraphRuns $gsiOutPara @(
    ,("text", "            Df  Sum Sq Mean Sq F value   Pr(>F)    ", "VerbatimChar")
    ,("br", $null, $null)
    ,("text", "Size_Class   1 0.13667 0.13667   63.31 1.57e-09 ***", "VerbatimChar")
    ,("br", $null, $null)
    ,("text", "Residuals   37 0.07988 0.00216                     ", "VerbatimChar")
    ,("br", $null, $null)
    ,("text", "---", "VerbatimChar")
    ,("br", $null, $null)
    ,("text", "Signif. codes:  0 '***' 0.001 '**' 0.01 '*' 0.05 '.' 0.1 ' ' 1", "VerbatimChar")
)

# ------------------------------------------------------------------
# 5. PCI aov() call paragraph
# ------------------------------------------------------------------
$pciCallPara = $d.Paragraphs.Item(20)
Set-ParagraphRuns $pciCallPara @(
    ,("text", "pci_anova ", "NormalTok")
    ,("text", "<-", "OtherTok")
    ,("text", " ", "NormalTok")
    ,("text", "aov", "FunctionTok")
    ,("text", "(PCI ", "NormalTok")
    ,("text", "~", "SpecialCharTok")
    ,("text", " Size_Class, ", "NormalTok")
    ,("text", "data =", "AttributeTok")
    ,("text", " cots_data)", "NormalTok")
    ,("br", $null, $null)
    ,("text", "summary", "FunctionTok")
    ,("text", "(pci_anova)", "NormalTok")
)

# ------------------------------------------------------------------
# 6. PCI summary() output paragraph
# ------------------------------------------------------------------
$pciOutPara = $d.Paragraphs.Item(21)
Set-ParagraphRuns $pciOutPara @(
    ,("text", "            Df Sum Sq Mean Sq F value   Pr(>F)    ", "VerbatimChar")
    ,("br", $null, $null)
    ,("text", "Size_Class   1  2.143  2.1428   63.31 1.57e-09 ***", "VerbatimChar")
    ,("br", $null, $null)
    ,("text", "Residuals   37  1.252  0.0338                     ", "VerbatimChar")
    ,("br", $null, $null)
    ,("text", "---", "VerbatimChar")
    ,("br", $null, $null)
    ,("text", "Signif. codes:  0 '***' 0.001 '**' 0.01 '*' 0.05 '.' 0.1 ' ' 1", "VerbatimChar")
)

Write-Host "Edits applied."
